$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 212.5
$ws.Range("J38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("N38").ClearContents()
$ws.Range("H44").Value = 74024
$ws.Range("J44").Value = 74024
$ws.Range("L44").Value = 74024
$ws.Range("N44").Value = -74948
$ws.Range("H58").Value = 149.41667
$ws.Range("I58").Value = 149.41667
$ws.Range("K58").Value = 448.25001
$ws.Range("M58").Value = -298.25001
$ws.Range("H112").Value = 1669334
$ws.Range("I112").Value = 699.5
$ws.Range("J112").Value = 2503651.2
$ws.Range("K112").Value = 2098.5
$ws.Range("L112").Value = 7510953.600000001
$ws.Range("M112").Value = -990.5
$ws.Range("N112").Value = -7513169.600000001
$ws.Range("H129").Value = 62503364
$ws.Range("I129").Value = 90909540
$ws.Range("K129").Value = 272728620
$ws.Range("M129").Value = -272723620
$ws.Range("H137").Value = 1661.9487
$ws.Range("I137").Value = 1753.7916
$ws.Range("J137").Value = 1515
$ws.Range("K137").Value = 5261.3748
$ws.Range("L137").Value = 4545
$ws.Range("M137").Value = -2711.3748
$ws.Range("N137").Value = -9645

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H26").Value = 3999.5
$ws.Range("I26").Value = 2999
$ws.Range("K26").Value = 2999
$ws.Range("M26").Value = -2669
$ws.Range("H32").Value = 7942.413
$ws.Range("I32").Value = 8313.584999999999
$ws.Range("J32").Value = 4898.8
$ws.Range("K32").Value = 8313.584999999999
$ws.Range("L32").Value = 4898.8
$ws.Range("M32").Value = -8026.584999999999
$ws.Range("N32").Value = -5472.8
$ws.Range("H41").Value = 1377
$ws.Range("I41").Value = 1377
$ws.Range("K41").Value = 1377
$ws.Range("M41").Value = -963
$ws.Range("H45").Value = 2856.7646
$ws.Range("I45").Value = 2417.4443
$ws.Range("K45").Value = 2417.4443
$ws.Range("M45").Value = -2040.4443
$ws.Range("H61").Value = 4499.0884
$ws.Range("I61").Value = 2896.442
$ws.Range("K61").Value = 2896.442
$ws.Range("M61").Value = -2684.442
$ws.Range("H82").Value = 110000
$ws.Range("J82").Value = 110000
$ws.Range("L82").Value = 110000
$ws.Range("N82").Value = -110722
$ws.Range("H85").Value = 110000
$ws.Range("J85").Value = 110000
$ws.Range("L85").Value = 110000
$ws.Range("N85").Value = -112496
$ws.Range("H136").Value = 4499.0884
$ws.Range("I136").Value = 2896.442
$ws.Range("K136").Value = 8689.326000000001
$ws.Range("M136").Value = -6139.326000000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H37").Value = 3553
$ws.Range("I37").Value = 348.66666
$ws.Range("J37").Value = 5475.6
$ws.Range("K37").Value = 348.66666
$ws.Range("L37").Value = 5475.6
$ws.Range("M37").Value = -211.66666
$ws.Range("N37").Value = -5749.6
$ws.Range("H86").Value = 8349809
$ws.Range("I86").Value = 13357954
$ws.Range("J86").Value = 2901.3333
$ws.Range("K86").Value = 13357954
$ws.Range("L86").Value = 2901.3333
$ws.Range("M86").Value = -13356831
$ws.Range("N86").Value = -5147.3333
$ws.Range("H89").Value = 8349809
$ws.Range("I89").Value = 13357954
$ws.Range("J89").Value = 2901.3333
$ws.Range("K89").Value = 66789770
$ws.Range("L89").Value = 14506.6665
$ws.Range("M89").Value = -66784154
$ws.Range("N89").Value = -25738.6665
$ws.Range("H105").Value = 4676
$ws.Range("I105").Value = 3105.2307
$ws.Range("K105").Value = 3105.2307
$ws.Range("M105").Value = -1358.2307
$ws.Range("H134").Value = 9310.718999999999
$ws.Range("I134").Value = 2647.6667
$ws.Range("J134").Value = 10000
$ws.Range("K134").Value = 7943.000100000001
$ws.Range("L134").Value = 30000
$ws.Range("M134").Value = -5408.000100000001
$ws.Range("N134").Value = -35070

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1964.7167
$ws.Range("I31").Value = 1152.275
$ws.Range("J31").Value = 3589.6
$ws.Range("K31").Value = 1152.275
$ws.Range("L31").Value = 3589.6
$ws.Range("M31").Value = -857.2750000000001
$ws.Range("N31").Value = -4179.6
$ws.Range("H34").Value = 1964.7167
$ws.Range("I34").Value = 1152.275
$ws.Range("J34").Value = 3589.6
$ws.Range("K34").Value = 1152.275
$ws.Range("L34").Value = 3589.6
$ws.Range("M34").Value = -950.2750000000001
$ws.Range("N34").Value = -3993.6
$ws.Range("H132").Value = 1484014.8
$ws.Range("I132").Value = 1669205
$ws.Range("J132").Value = 2493
$ws.Range("K132").Value = 5007615
$ws.Range("L132").Value = 7479
$ws.Range("M132").Value = -5005085
$ws.Range("N132").Value = -12539

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H22").Value = 2684.1667
$ws.Range("I22").Value = 221
$ws.Range("J22").Value = 15000
$ws.Range("K22").Value = 663
$ws.Range("L22").Value = 45000
$ws.Range("M22").Value = -494
$ws.Range("N22").Value = -45338
$ws.Range("H27").Value = 2684.1667
$ws.Range("I27").Value = 221
$ws.Range("J27").Value = 15000
$ws.Range("K27").Value = 663
$ws.Range("L27").Value = 45000
$ws.Range("M27").Value = -561
$ws.Range("N27").Value = -45204
$ws.Range("H47").Value = 3623.75
$ws.Range("J47").Value = 747.5
$ws.Range("L47").Value = 2242.5
$ws.Range("N47").Value = -3104.5
$ws.Range("H129").Value = 1091.5
$ws.Range("I129").Value = 388.125
$ws.Range("J129").Value = 2967.1667
$ws.Range("K129").Value = 1164.375
$ws.Range("L129").Value = 8901.500100000001
$ws.Range("M129").Value = 3835.625
$ws.Range("N129").Value = -18901.5001
$ws.Range("H131").Value = 1743.2
$ws.Range("J131").Value = 2144.4167
$ws.Range("L131").Value = 6433.250100000001
$ws.Range("N131").Value = -16513.2501
$ws.Range("H137").Value = 3577.4167
$ws.Range("J137").Value = 3444.3333
$ws.Range("L137").Value = 10332.9999
$ws.Range("N137").Value = -20532.9999
$ws.Range("H141").Value = 18707.143
$ws.Range("I141").Value = 16158.333
$ws.Range("K141").Value = 48474.999
$ws.Range("M141").Value = -43294.999

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7775.909
$ws.Range("I70").Value = 7089.3335
$ws.Range("K70").Value = 7089.3335
$ws.Range("M70").Value = -6819.3335
$ws.Range("H73").Value = 7775.909
$ws.Range("I73").Value = 7089.3335
$ws.Range("K73").Value = 7089.3335
$ws.Range("M73").Value = -6153.3335
$ws.Range("H80").Value = 3089.3
$ws.Range("I80").Value = 2880.7273
$ws.Range("J80").Value = 3662.875
$ws.Range("K80").Value = 2880.7273
$ws.Range("L80").Value = 3662.875
$ws.Range("M80").Value = -1882.7273
$ws.Range("N80").Value = -5658.875
$ws.Range("H83").Value = 3089.3
$ws.Range("I83").Value = 2880.7273
$ws.Range("J83").Value = 3662.875
$ws.Range("K83").Value = 14403.6365
$ws.Range("L83").Value = 18314.375
$ws.Range("M83").Value = -9411.636500000001
$ws.Range("N83").Value = -28298.375

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6118.3335
$ws.Range("I7").Value = 4482
$ws.Range("J7").Value = 7550.125
$ws.Range("K7").Value = 4482
$ws.Range("L7").Value = 7550.125
$ws.Range("M7").Value = -4370
$ws.Range("N7").Value = -7774.125
$ws.Range("H22").Value = 3498.125
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H27").Value = 3498.125
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("H46").Value = 2972.6667
$ws.Range("I46").Value = 2666.6667
$ws.Range("J46").Value = 3156.2666
$ws.Range("K46").Value = 2666.6667
$ws.Range("L46").Value = 3156.2666
$ws.Range("M46").Value = -2478.6667
$ws.Range("N46").Value = -3532.2666
$ws.Range("H93").Value = 2905
$ws.Range("I93").Value = 2249
$ws.Range("K93").Value = 2249
$ws.Range("M93").Value = -1001
$ws.Range("H126").Value = 6118.3335
$ws.Range("I126").Value = 4482
$ws.Range("J126").Value = 7550.125
$ws.Range("K126").Value = 13446
$ws.Range("L126").Value = 22650.375
$ws.Range("M126").Value = -10976
$ws.Range("N126").Value = -27590.375
$ws.Range("H132").Value = 3081.468
$ws.Range("I132").Value = 3060.658
$ws.Range("K132").Value = 9181.974
$ws.Range("M132").Value = -6651.974

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1299
$ws.Range("I107").Value = 600
$ws.Range("K107").Value = 1800
$ws.Range("M107").Value = 120
$ws.Range("H132").Value = 5316.8
$ws.Range("I132").Value = 6537.5713
$ws.Range("J132").Value = 2468.3333
$ws.Range("K132").Value = 19612.7139
$ws.Range("L132").Value = 7404.999899999999
$ws.Range("M132").Value = -17082.7139
$ws.Range("N132").Value = -12464.9999

Write-Output "Applied all Faerie_Profits updates"